$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New headers for the velocity columns ---
$ws.Range("D1").Value = "vx"
$ws.Range("E1").Value = "vy"

# --- Velocity data for b1 (row2), b2 (row3), b3 (row4), b4 (row5) ---
$ws.Range("D2").Value = 4
$ws.Range("E2").Value = 2
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = 5
$ws.Range("D4").Value = 5
$ws.Range("E4").Value = 4
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 0

# --- Scratch RANDBETWEEN helper column (rows 2-8) ---
$ws.Range("F2").Formula = "=RANDBETWEEN(0, 5)"
$ws.Range("F3").Formula = "=RANDBETWEEN(0, 5)"
$ws.Range("F4").Formula = "=RANDBETWEEN(0, 5)"
$ws.Range("F5").Formula = "=RANDBETWEEN(0, 5)"
$ws.Range("F6").Formula = "=RANDBETWEEN(0, 5)"
$ws.Range("F7").Formula = "=RANDBETWEEN(0, 5)"
$ws.Range("F8").Formula = "=RANDBETWEEN(0, 5)"

# --- Average velocity excluding b2, used for the velocity matching rule ---
$ws.Range("D7").Value = "avg velocity excluding b2"
$ws.Range("D8").Formula = "=(D2+D4+D5)/3"
$ws.Range("E8").Formula = "=(E2+E4+E5)/3"

# --- Move/resize the scatter chart to make room for the new columns ---
$co = $ws.ChartObjects().Item(1)
$co.Left = 470.0625
$co.Top = 337
$co.Width = 277.1875
$co.Height = 216

# --- Update the active selection ---
$ws.Range("E9").Select() | Out-Null
